$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in B27 (week 25-01/05/2016) to add an extra 1 hour entry.
$ws.Range("B27").Formula = "= 4.5 + 4 + 1.5 + 1 + 2.5 + 1 + 1.5 + 2.5 + 3.25 + 1.5 + 2.5 + 1"

# Move/update the active selection to B28, matching the saved sheet view state.
$ws.Range("B28").Select()
